# Updates cryptos list prices/volumes (Price = column D, Volume(1h) = column E).
# Column D values are stored as text that merely look numeric (e.g. "63.250.02",
# "0.0539", "18.80"); a leading apostrophe forces Excel to keep them as literal
# text instead of coercing them to numbers (which would also drop meaningful
# trailing/leading zeros). In a PowerShell single-quoted string literal, an
# embedded apostrophe must be written as '' (doubled), so a value like
# "'63.250.02" is written below as '''63.250.02'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.250.02'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '''2.455.79'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''572.43'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '''147.55'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '''0.540'
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('D9').Value = '''2.449.75'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '''0.111'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').Value = '''5.26'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '''0.355'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').Value = '''27.24'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('D15').Value = '''0.0000180'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').Value = '''2.906.09'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '''63.141.15'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '''2.466.56'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''11.37'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').Value = '''7.39'
$ws.Range('E20').Value = '  +6.66%  '
$ws.Range('D21').Value = '''328.46'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').Value = '''4.20'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +13.90%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''65.46'
$ws.Range('E25').Value = '  -2.74%  '
$ws.Range('D26').Value = '''623.58'
$ws.Range('E26').Value = '  +8.11%  '
$ws.Range('D27').Value = '''8.98'
$ws.Range('E27').Value = '  +3.17%  '
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('D29').Value = '''2.605.14'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '''1.51'
$ws.Range('E30').Value = '  +4.50%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('D32').Value = '''8.30'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('D34').Value = '''0.141'
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('D35').Value = '''5.21'
$ws.Range('E35').Value = '  +7.49%  '
$ws.Range('D36').Value = '''1.53'
$ws.Range('E36').Value = '  -1.59%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '''0.381'
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('D39').Value = '''5.48'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = '''18.80'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = '''2.72'
$ws.Range('E41').Value = '  +12.14%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '''145.83'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.79'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').Value = '''149.55'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '''3.77'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('D47').Value = '''21.34'
$ws.Range('E47').Value = '  +4.02%  '
$ws.Range('D48').Value = '''0.0539'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').Value = '''0.602'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '''0.0235'
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('D51').Value = '''0.0917'
$ws.Range('E51').Value = '  -0.81%  '
